$p = $ppt.ActivePresentation

# Delete slides 3,4,5,6 (old "Title/Content" filler slides).
# Delete from the end backward so indices stay valid.
$p.Slides.Item(6).Delete()
$p.Slides.Item(5).Delete()
$p.Slides.Item(4).Delete()
$p.Slides.Item(3).Delete()

Write-Output ("Slide count after: " + $p.Slides.Count)
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    Write-Output ("Slide " + $i + ": " + $s.Shapes.Count + " shapes")
}
